$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B6 value from 20.0 to 50.0
$ws.Range("B6").Value = 50.0

# Add new cell F6 with value 240.0, matching the format/style used by F7 (s=5)
$ws.Range("F6").Value = 240.0
$ws.Range("F7").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# Update G6, H6, I6 values
$ws.Range("G6").Value = -1.031009
$ws.Range("H6").Value = -1.581433
$ws.Range("I6").Value = -3.491137
